# Applies the Ramuh_Profits.xlsx edits (scheduled-runner price/profit refresh)
# across the ALC, ARM, CRP, CUL, GSM, LTW, and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value2 = 311.66666
$ws.Range("I61").Value2 = 311.66666
$ws.Range("K61").Value2 = 934.9999799999999
$ws.Range("M61").Value2 = -762.9999799999999

$ws.Range("H62").Value2 = 1881.1875
$ws.Range("I62").Value2 = 2026.25
$ws.Range("J62").Value2 = 1832.8334
$ws.Range("K62").Value2 = 2026.25
$ws.Range("L62").Value2 = 1832.8334
$ws.Range("M62").Value2 = -1402.25
$ws.Range("N62").Value2 = -3080.8334

$ws.Range("H64").Value2 = 34485560
$ws.Range("I64").Value2 = 50002280
$ws.Range("J64").Value2 = 3962.2222
$ws.Range("K64").Value2 = 50002280
$ws.Range("L64").Value2 = 3962.2222
$ws.Range("M64").Value2 = -50002032
$ws.Range("N64").Value2 = -4458.2222

$ws.Range("H65").Value2 = 1881.1875
$ws.Range("I65").Value2 = 2026.25
$ws.Range("J65").Value2 = 1832.8334
$ws.Range("K65").Value2 = 10131.25
$ws.Range("L65").Value2 = 9164.166999999999
$ws.Range("M65").Value2 = -7011.25
$ws.Range("N65").Value2 = -15404.167

$ws.Range("H67").Value2 = 34485560
$ws.Range("I67").Value2 = 50002280
$ws.Range("J67").Value2 = 3962.2222
$ws.Range("K67").Value2 = 50002280
$ws.Range("L67").Value2 = 3962.2222
$ws.Range("M67").Value2 = -50001422
$ws.Range("N67").Value2 = -5678.2222

$ws.Range("H82").Value2 = 372.5
$ws.Range("I82").Value2 = 372.5
$ws.Range("K82").Value2 = 1117.5
$ws.Range("M82").Value2 = -711.5

$ws.Range("H85").Value2 = 372.5
$ws.Range("I85").Value2 = 372.5
$ws.Range("K85").Value2 = 1117.5
$ws.Range("M85").Value2 = 286.5

$ws.Range("H92").Value2 = 750.8
$ws.Range("I92").Value2 = 362.66666
$ws.Range("J92").Value2 = 2303.3333
$ws.Range("K92").Value2 = 362.66666
$ws.Range("L92").Value2 = 2303.3333
$ws.Range("M92").Value2 = 885.33334
$ws.Range("N92").Value2 = -4799.3333

$ws.Range("H93").Value2 = 35000.668
$ws.Range("J93").Value2 = 35001
$ws.Range("L93").Value2 = 35001
$ws.Range("N93").Value2 = -39993

$ws.Range("H98").Value2 = 1758.9354
$ws.Range("I98").Value2 = 1398.7142
$ws.Range("J98").Value2 = 2515.4
$ws.Range("K98").Value2 = 1398.7142
$ws.Range("L98").Value2 = 2515.4
$ws.Range("M98").Value2 = 99.28580000000011
$ws.Range("N98").Value2 = -5511.4

$ws.Range("H100").Value2 = 1531.4286
$ws.Range("I100").Value2 = 1531.4286
$ws.Range("J100").Value2 = 0
$ws.Range("K100").Value2 = 1531.4286
$ws.Range("L100").Value2 = 0
$ws.Range("M100").Value2 = -990.4286
$ws.Range("N100").ClearContents()

$ws.Range("H107").Value2 = 496.64285
$ws.Range("I107").Value2 = 504.07693
$ws.Range("K107").Value2 = 504.07693
$ws.Range("M107").Value2 = 1415.92307

$ws.Range("H109").Value2 = 48000
$ws.Range("J109").Value2 = 48000
$ws.Range("L109").Value2 = 48000
$ws.Range("N109").Value2 = -50774

$ws.Range("H118").Value2 = 1806.1428
$ws.Range("I118").Value2 = 626.5714
$ws.Range("J118").Value2 = 2985.7144
$ws.Range("K118").Value2 = 1879.7142
$ws.Range("L118").Value2 = 8957.143199999999
$ws.Range("M118").Value2 = -222.7142000000001
$ws.Range("N118").Value2 = -12271.1432

$ws.Range("H122").Value2 = 1758.9354
$ws.Range("I122").Value2 = 1398.7142
$ws.Range("J122").Value2 = 2515.4
$ws.Range("K122").Value2 = 4196.142599999999
$ws.Range("L122").Value2 = 7546.200000000001
$ws.Range("M122").Value2 = -1746.142599999999
$ws.Range("N122").Value2 = -12446.2

$ws.Range("H123").Value2 = 27975
$ws.Range("J123").Value2 = 27975
$ws.Range("L123").Value2 = 27975
$ws.Range("N123").Value2 = -37775

$ws.Range("H125").Value2 = 11746.546
$ws.Range("I125").Value2 = 40704
$ws.Range("J125").Value2 = 887.5
$ws.Range("K125").Value2 = 366336
$ws.Range("L125").Value2 = 7987.5
$ws.Range("M125").Value2 = -363876
$ws.Range("N125").Value2 = -12907.5

$ws.Range("H127").Value2 = 860.0741
$ws.Range("I127").Value2 = 713.2941
$ws.Range("K127").Value2 = 2139.8823
$ws.Range("M127").Value2 = 2820.1177

$ws.Range("H129").Value2 = 917.1070999999999
$ws.Range("I129").Value2 = 894.25
$ws.Range("J129").Value2 = 920.9167
$ws.Range("K129").Value2 = 2682.75
$ws.Range("L129").Value2 = 2762.7501
$ws.Range("M129").Value2 = 2317.25
$ws.Range("N129").Value2 = -12762.7501

$ws.Range("H132").Value2 = 1907.0358
$ws.Range("I132").Value2 = 1651.8
$ws.Range("J132").Value2 = 4034
$ws.Range("K132").Value2 = 4955.4
$ws.Range("L132").Value2 = 12102
$ws.Range("M132").Value2 = -2425.4
$ws.Range("N132").Value2 = -17162

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value2 = 76166.664
$ws.Range("J69").Value2 = 76166.664
$ws.Range("L69").Value2 = 76166.664
$ws.Range("N69").Value2 = -77664.664

$ws.Range("H72").Value2 = 76166.664
$ws.Range("J72").Value2 = 76166.664
$ws.Range("L72").Value2 = 228499.992
$ws.Range("N72").Value2 = -235987.992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1810.79
$ws.Range("I31").Value2 = 1121.82
$ws.Range("J31").Value2 = 2499.76
$ws.Range("K31").Value2 = 1121.82
$ws.Range("L31").Value2 = 2499.76
$ws.Range("M31").Value2 = -826.8199999999999
$ws.Range("N31").Value2 = -3089.76

$ws.Range("H34").Value2 = 1810.79
$ws.Range("I34").Value2 = 1121.82
$ws.Range("J34").Value2 = 2499.76
$ws.Range("K34").Value2 = 1121.82
$ws.Range("L34").Value2 = 2499.76
$ws.Range("M34").Value2 = -919.8199999999999
$ws.Range("N34").Value2 = -2903.76

$ws.Range("H99").Value2 = 1478.5
$ws.Range("I99").Value2 = 1504.8
$ws.Range("J99").Value2 = 1434.6666
$ws.Range("K99").Value2 = 1504.8
$ws.Range("L99").Value2 = 1434.6666
$ws.Range("M99").Value2 = -6.799999999999955
$ws.Range("N99").Value2 = -4430.6666

$ws.Range("H107").Value2 = 680.0833
$ws.Range("I107").Value2 = 684.0625
$ws.Range("J107").Value2 = 648.25
$ws.Range("K107").Value2 = 684.0625
$ws.Range("L107").Value2 = 648.25
$ws.Range("M107").Value2 = 1235.9375
$ws.Range("N107").Value2 = -4488.25

$ws.Range("H126").Value2 = 1478.5
$ws.Range("I126").Value2 = 1504.8
$ws.Range("J126").Value2 = 1434.6666
$ws.Range("K126").Value2 = 4514.4
$ws.Range("L126").Value2 = 4303.9998
$ws.Range("M126").Value2 = -2044.4
$ws.Range("N126").Value2 = -9243.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value2 = 1005122.9
$ws.Range("I133").Value2 = 5404.2856
$ws.Range("J133").Value2 = 3337799.8
$ws.Range("K133").Value2 = 16212.8568
$ws.Range("L133").Value2 = 10013399.4
$ws.Range("M133").Value2 = -11152.8568
$ws.Range("N133").Value2 = -10023519.4

$ws.Range("H137").Value2 = 3648488.5
$ws.Range("I137").Value2 = 111846
$ws.Range("K137").Value2 = 335538
$ws.Range("M137").Value2 = -330438

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 5475
$ws.Range("J122").Value2 = 7033.3335
$ws.Range("L122").Value2 = 21100.0005
$ws.Range("N122").Value2 = -26000.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 3264.25
$ws.Range("I7").Value2 = 3420.8
$ws.Range("J7").Value2 = 3003.3333
$ws.Range("K7").Value2 = 3420.8
$ws.Range("L7").Value2 = 3003.3333
$ws.Range("M7").Value2 = -3308.8
$ws.Range("N7").Value2 = -3227.3333

$ws.Range("H40").Value2 = 2380.5454
$ws.Range("I40").Value2 = 2242.3333
$ws.Range("J40").Value2 = 3002.5
$ws.Range("K40").Value2 = 2242.3333
$ws.Range("L40").Value2 = 3002.5
$ws.Range("M40").Value2 = -2106.3333
$ws.Range("N40").Value2 = -3274.5

$ws.Range("H48").Value2 = 5566.6665
$ws.Range("I48").Value2 = 3850
$ws.Range("K48").Value2 = 3850
$ws.Range("M48").Value2 = -3189

$ws.Range("H122").Value2 = 80831.84
$ws.Range("I122").Value2 = 1000004
$ws.Range("J122").Value2 = 4234.1665
$ws.Range("K122").Value2 = 3000012
$ws.Range("L122").Value2 = 12702.4995
$ws.Range("M122").Value2 = -2997562
$ws.Range("N122").Value2 = -17602.4995

$ws.Range("H126").Value2 = 3264.25
$ws.Range("I126").Value2 = 3420.8
$ws.Range("J126").Value2 = 3003.3333
$ws.Range("K126").Value2 = 10262.4
$ws.Range("L126").Value2 = 9009.999899999999
$ws.Range("M126").Value2 = -7792.400000000001
$ws.Range("N126").Value2 = -13949.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value2 = 1096.9333
$ws.Range("I126").Value2 = 932.6667
$ws.Range("J126").Value2 = 1343.3334
$ws.Range("K126").Value2 = 2798.0001
$ws.Range("L126").Value2 = 4030.0002
$ws.Range("M126").Value2 = -328.0001000000002
$ws.Range("N126").Value2 = -8970.0002
